$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0531941391419865
$ws.Range("C2").Value = 0.0320326737718645
$ws.Range("D2").Value = 0.0743556045121085

$ws.Range("B3").Value = 0.0575653552523801
$ws.Range("C3").Value = 0.0288289010822672
$ws.Range("D3").Value = 0.086301809422493

$ws.Range("B4").Value = 0.0207598905329338
$ws.Range("D4").Value = 0.0876379770828152

$ws.Range("B5").Value = 0.091733891587786
$ws.Range("C5").Value = 0.0297806852548251
$ws.Range("D5").Value = 0.153687097920747

$ws.Range("B6").Value = 0.0838719260554863
$ws.Range("C6").Value = 0.00843318620479258

$ws.Range("B8").Value = 0.0769389816500664
$ws.Range("C8").Value = -0.00563473232866796

$ws.Range("B9").Value = 0.030970014094514
$ws.Range("D9").Value = 0.0989919038236001

$ws.Range("B10").Value = -0.00623267936164858

$ws.Range("B11").Value = 0.0342307930815836
$ws.Range("C11").Value = -0.000817353665763139
$ws.Range("D11").Value = 0.0692789398289303

$ws.Range("B12").Value = 0.0740914674696293
$ws.Range("C12").Value = 0.015039994457167
$ws.Range("D12").Value = 0.133142940482092
